# Weekly refresh: insert a new "Cilantro" price record at the top of the
# data block (row 9) and push the existing rows down by one, so the most
# recent observation (2023-05-05) leads the series and the oldest row
# (previously row 87) now lands on row 88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:87 down one row, creating a blank row 9.
$ws.Rows("9:9").Insert(-4121)

# Populate the new row 9 with this week's record.
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(9, 3).Value = "Maule"
$ws.Cells.Item(9, 4).Value = 45051
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 100112040
$ws.Cells.Item(9, 7).Value = "Cilantro"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 150
$ws.Cells.Item(9, 11).Value = 7000
$ws.Cells.Item(9, 12).Value = 7000
$ws.Cells.Item(9, 13).Value = 7000
$ws.Cells.Item(9, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 194
$ws.Cells.Item(9, 17).Value = 36
$ws.Cells.Item(9, 18).Value = "Hortaliza"
